# Add a new "priority" column (I) to the exchanges sheet.
#
# A small, hand-picked subset of exchanges gets a numeric priority rank
# (used for sorting/searching "more generally" per the commit message);
# everything else in column I stays blank, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the other header cells in row 1 (C1:H1).
$ws.Range("I1").Value = "priority"
$ws.Range("C1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Exchange code (column A) -> priority value (column I), by row.
$priorities = @{
    4  = 1    # AX - ASX - ALL MARKETS
    12 = 7    # CN - CANADIAN NATIONAL STOCK EXCHANGE
    13 = 21   # CO - OMX NORDIC EXCHANGE COPENHAGEN A/S
    15 = 6    # DB - DUBAI FINANCIAL MARKET
    18 = 18   # F  - DEUTSCHE BOERSE AG
    19 = 20   # HE - NASDAQ OMX HELSINKI LTD.
    20 = 5    # HK - HONG KONG EXCHANGES AND CLEARING LTD
    22 = 14   # IC - NASDAQ OMX ICELAND
    28 = 17   # KQ - KOREA EXCHANGE (KOSDAQ)
    29 = 16   # KS - KOREA EXCHANGE (STOCK MARKET)
    30 = 4    # L  - LONDON STOCK EXCHANGE
    39 = 3    # NZ - NEW ZEALAND EXCHANGE LTD
    42 = 15   # PR - PRAGUE STOCK EXCHANGE
    51 = 8    # ST - NASDAQ OMX NORDIC
    52 = 9    # SW - SWISS EXCHANGE
    53 = 10   # SZ - SHENZHEN STOCK EXCHANGE
    54 = 11   # T  - TOKYO STOCK EXCHANGE-TOKYO PRO MARKET
    57 = 12   # TO - TORONTO STOCK EXCHANGE
    59 = 2    # US - US exchanges (NYSE, Nasdaq)
    61 = 11   # VI - Vienna Stock Exchange
    64 = 13   # WA - WARSAW STOCK EXCHANGE/EQUITIES/MAIN MARKET
}

foreach ($row in $priorities.Keys) {
    $ws.Cells.Item($row, 9).Value = $priorities[$row]
}
